$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5293595
$ws.Range("I51").Value = 3666.3333
$ws.Range("J51").Value = 6175249.5
$ws.Range("K51").Value = 3666.3333
$ws.Range("L51").Value = 6175249.5
$ws.Range("M51").Value = -3182.3333
$ws.Range("N51").Value = -6176217.5

$ws.Range("H58").Value = 1250
$ws.Range("I58").Value = 284.22223
$ws.Range("J58").Value = 2215.7778
$ws.Range("K58").Value = 852.66669
$ws.Range("L58").Value = 6647.3334
$ws.Range("M58").Value = -702.66669
$ws.Range("N58").Value = -6947.3334

$ws.Range("H86").Value = 249438.67
$ws.Range("I86").Value = 446888.8
$ws.Range("J86").Value = 2626
$ws.Range("K86").Value = 446888.8
$ws.Range("L86").Value = 2626
$ws.Range("M86").Value = -445765.8
$ws.Range("N86").Value = -4872

$ws.Range("H89").Value = 249438.67
$ws.Range("I89").Value = 446888.8
$ws.Range("J89").Value = 2626
$ws.Range("K89").Value = 2234444
$ws.Range("L89").Value = 13130
$ws.Range("M89").Value = -2228828
$ws.Range("N89").Value = -24362

$ws.Range("H116").Value = 1703015.6
$ws.Range("I116").Value = 2034757.2
$ws.Range("J116").Value = 2839.625
$ws.Range("K116").Value = 2034757.2
$ws.Range("L116").Value = 2839.625
$ws.Range("M116").Value = -2031315.2
$ws.Range("N116").Value = -9723.625

$ws.Range("H121").Value = 1528.3334
$ws.Range("I121").Value = 5095
$ws.Range("J121").Value = 1082.5
$ws.Range("K121").Value = 15285
$ws.Range("L121").Value = 3247.5
$ws.Range("M121").Value = -13538
$ws.Range("N121").Value = -6741.5

$ws.Range("H132").Value = 2824.459
$ws.Range("I132").Value = 2619.84
$ws.Range("J132").Value = 3754.5454
$ws.Range("K132").Value = 7859.52
$ws.Range("L132").Value = 11263.6362
$ws.Range("M132").Value = -5329.52
$ws.Range("N132").Value = -16323.6362

$ws.Range("H138").Value = 2787.4219
$ws.Range("I138").Value = 1674.6923
$ws.Range("J138").Value = 3071.0588
$ws.Range("K138").Value = 5024.0769
$ws.Range("L138").Value = 9213.1764
$ws.Range("M138").Value = 115.9231
$ws.Range("N138").Value = -19493.1764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 694.8
$ws.Range("I2").Value = 663.5333000000001
$ws.Range("J2").Value = 788.6
$ws.Range("K2").Value = 663.5333000000001
$ws.Range("L2").Value = 788.6
$ws.Range("M2").Value = -550.5333000000001
$ws.Range("N2").Value = -1014.6

$ws.Range("H32").Value = 16460.055
$ws.Range("I32").Value = 5460.5625
$ws.Range("J32").Value = 28738.559
$ws.Range("K32").Value = 5460.5625
$ws.Range("L32").Value = 28738.559
$ws.Range("M32").Value = -5173.5625
$ws.Range("N32").Value = -29312.559

$ws.Range("H45").Value = 2140.5789
$ws.Range("I45").Value = 2243.923
$ws.Range("J45").Value = 1916.6666
$ws.Range("K45").Value = 2243.923
$ws.Range("L45").Value = 1916.6666
$ws.Range("M45").Value = -1866.923
$ws.Range("N45").Value = -2670.6666

$ws.Range("H69").Value = 79750
$ws.Range("J69").Value = 79750
$ws.Range("L69").Value = 79750
$ws.Range("N69").Value = -81248

$ws.Range("H72").Value = 79750
$ws.Range("J72").Value = 79750
$ws.Range("L72").Value = 239250
$ws.Range("N72").Value = -246738

$ws.Range("H92").Value = 32000
$ws.Range("J92").Value = 32000
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992

$ws.Range("H116").Value = 694.8
$ws.Range("I116").Value = 663.5333000000001
$ws.Range("J116").Value = 788.6
$ws.Range("K116").Value = 663.5333000000001
$ws.Range("L116").Value = 788.6
$ws.Range("M116").Value = 1630.4667
$ws.Range("N116").Value = -5376.6

$ws.Range("H122").Value = 1716.2122
$ws.Range("I122").Value = 1382.8077
$ws.Range("K122").Value = 4148.4231
$ws.Range("M122").Value = -1698.4231

$ws.Range("H132").Value = 4104.1665
$ws.Range("I132").Value = 3552.4
$ws.Range("J132").Value = 4498.2856
$ws.Range("K132").Value = 10657.2
$ws.Range("L132").Value = 13494.8568
$ws.Range("M132").Value = -8127.200000000001
$ws.Range("N132").Value = -18554.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 694.8
$ws.Range("I3").Value = 663.5333000000001
$ws.Range("J3").Value = 788.6
$ws.Range("K3").Value = 663.5333000000001
$ws.Range("L3").Value = 788.6
$ws.Range("M3").Value = -549.5333000000001
$ws.Range("N3").Value = -1016.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 3000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H107").Value = 1630.0834
$ws.Range("I107").Value = 2169.6
$ws.Range("J107").Value = 1244.7142
$ws.Range("K107").Value = 2169.6
$ws.Range("L107").Value = 1244.7142
$ws.Range("M107").Value = -249.5999999999999
$ws.Range("N107").Value = -5084.7142

$ws.Range("H132").Value = 2764.7144
$ws.Range("I132").Value = 1736.75
$ws.Range("J132").Value = 4135.3335
$ws.Range("K132").Value = 5210.25
$ws.Range("L132").Value = 12406.0005
$ws.Range("M132").Value = -2680.25
$ws.Range("N132").Value = -17466.0005

$ws.Range("H134").Value = 4632.231
$ws.Range("I134").Value = 5824.2144
$ws.Range("J134").Value = 3241.5833
$ws.Range("K134").Value = 17472.6432
$ws.Range("L134").Value = 9724.749899999999
$ws.Range("M134").Value = -14937.6432
$ws.Range("N134").Value = -14794.7499

$ws.Range("H141").Value = 66320.47
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 66320.47
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 66320.47
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -76680.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 100
$ws.Range("J23").Value = 100
$ws.Range("L23").Value = 300
$ws.Range("N23").Value = -770

$ws.Range("H92").Value = 428.125
$ws.Range("I92").Value = 412.5
$ws.Range("J92").Value = 475
$ws.Range("K92").Value = 1237.5
$ws.Range("L92").Value = 1425
$ws.Range("M92").Value = 10.5
$ws.Range("N92").Value = -3921

$ws.Range("H113").Value = 922.2963
$ws.Range("I113").Value = 448.2857
$ws.Range("J113").Value = 1088.2
$ws.Range("K113").Value = 1344.8571
$ws.Range("L113").Value = 3264.6
$ws.Range("M113").Value = 825.1428999999998
$ws.Range("N113").Value = -7604.6

$ws.Range("H122").Value = 736
$ws.Range("J122").Value = 1186.4445
$ws.Range("L122").Value = 10678.0005
$ws.Range("N122").Value = -15578.0005

$ws.Range("H131").Value = 943.34
$ws.Range("J131").Value = 1057.175
$ws.Range("L131").Value = 3171.525
$ws.Range("N131").Value = -13251.525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2104.6
$ws.Range("J80").Value = 2055.4055
$ws.Range("L80").Value = 2055.4055
$ws.Range("N80").Value = -4051.4055

$ws.Range("H83").Value = 2104.6
$ws.Range("J83").Value = 2055.4055
$ws.Range("L83").Value = 10277.0275
$ws.Range("N83").Value = -20261.0275

$ws.Range("H102").Value = 880.55554
$ws.Range("I102").Value = 849.06665
$ws.Range("K102").Value = 849.06665
$ws.Range("M102").Value = 772.93335

$ws.Range("H113").Value = 1375.1333
$ws.Range("I113").Value = 1011.5455
$ws.Range("J113").Value = 2375
$ws.Range("K113").Value = 1011.5455
$ws.Range("L113").Value = 2375
$ws.Range("M113").Value = 1158.4545
$ws.Range("N113").Value = -6715

$ws.Range("H132").Value = 2768.3333
$ws.Range("I132").Value = 1644.6
$ws.Range("J132").Value = 3571
$ws.Range("K132").Value = 4933.799999999999
$ws.Range("L132").Value = 10713
$ws.Range("M132").Value = -2403.799999999999
$ws.Range("N132").Value = -15773

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1638.4286
$ws.Range("I7").Value = 1066.75
$ws.Range("J7").Value = 1936.6957
$ws.Range("K7").Value = 1066.75
$ws.Range("L7").Value = 1936.6957
$ws.Range("M7").Value = -954.75
$ws.Range("N7").Value = -2160.6957

$ws.Range("H22").Value = 998.3333
$ws.Range("I22").Value = 497.14285
$ws.Range("J22").Value = 1700
$ws.Range("K22").Value = 497.14285
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = -202.14285
$ws.Range("N22").Value = -2290

$ws.Range("H27").Value = 998.3333
$ws.Range("I27").Value = 497.14285
$ws.Range("J27").Value = 1700
$ws.Range("K27").Value = 497.14285
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = -390.14285
$ws.Range("N27").Value = -1914

$ws.Range("H61").Value = 3827.625
$ws.Range("I61").Value = 3786
$ws.Range("J61").Value = 3952.5
$ws.Range("K61").Value = 3786
$ws.Range("L61").Value = 3952.5
$ws.Range("M61").Value = -3584
$ws.Range("N61").Value = -4356.5

$ws.Range("H93").Value = 1918.75
$ws.Range("I93").Value = 1535.5555
$ws.Range("K93").Value = 1535.5555
$ws.Range("M93").Value = -287.5554999999999

$ws.Range("H113").Value = 3827.625
$ws.Range("I113").Value = 3786
$ws.Range("J113").Value = 3952.5
$ws.Range("K113").Value = 3786
$ws.Range("L113").Value = 3952.5
$ws.Range("M113").Value = -1616
$ws.Range("N113").Value = -8292.5

$ws.Range("H126").Value = 1638.4286
$ws.Range("I126").Value = 1066.75
$ws.Range("J126").Value = 1936.6957
$ws.Range("K126").Value = 3200.25
$ws.Range("L126").Value = 5810.0871
$ws.Range("M126").Value = -730.25
$ws.Range("N126").Value = -10750.0871

$ws.Range("H132").Value = 50024660
$ws.Range("J132").Value = 5100.6665
$ws.Range("L132").Value = 15301.9995
$ws.Range("N132").Value = -20361.9995

$ws.Range("H134").Value = 49953.848
$ws.Range("J134").Value = 49953.848
$ws.Range("L134").Value = 49953.848
$ws.Range("N134").Value = -60093.848

$ws.Range("H135").Value = 47107
$ws.Range("J135").Value = 47107
$ws.Range("L135").Value = 47107
$ws.Range("N135").Value = -57247

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3247.5881
$ws.Range("I126").Value = 4282.1816
$ws.Range("J126").Value = 1350.8334
$ws.Range("K126").Value = 12846.5448
$ws.Range("L126").Value = 4052.5002
$ws.Range("M126").Value = -10376.5448
$ws.Range("N126").Value = -8992.5002

$ws.Range("H132").Value = 1779.9
$ws.Range("I132").Value = 753.3
$ws.Range("J132").Value = 2293.2
$ws.Range("K132").Value = 2259.9
$ws.Range("L132").Value = 6879.599999999999
$ws.Range("M132").Value = 270.1000000000004
$ws.Range("N132").Value = -11939.6
